$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1001.11536
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1001.11536
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3003.34608
$ws.Range("N17").Value = -3339.34608

$ws.Range("H33").Value = 778.1111
$ws.Range("I33").Value = 966.8333
$ws.Range("J33").Value = 400.66666
$ws.Range("K33").Value = 966.8333
$ws.Range("L33").Value = 400.66666
$ws.Range("M33").Value = -737.8333
$ws.Range("N33").Value = -858.66666

$ws.Range("H80").Value = 1214.4445
$ws.Range("I80").Value = 718.5714
$ws.Range("J80").Value = 2950
$ws.Range("K80").Value = 2155.7142
$ws.Range("L80").Value = 8850
$ws.Range("M80").Value = -1157.7142
$ws.Range("N80").Value = -10846

$ws.Range("H83").Value = 1214.4445
$ws.Range("I83").Value = 718.5714
$ws.Range("J83").Value = 2950
$ws.Range("K83").Value = 6467.1426
$ws.Range("L83").Value = 26550
$ws.Range("M83").Value = -1475.1426
$ws.Range("N83").Value = -36534

$ws.Range("H106").Value = 4256.6113
$ws.Range("I106").Value = 4256.6113
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4256.6113
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3625.6113

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 0

$ws.Range("H115").Value = 1384.2
$ws.Range("I115").Value = 1384.2
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 4152.6
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -2585.6

$ws.Range("H116").Value = 22235204
$ws.Range("I116").Value = 37055372
$ws.Range("J116").Value = 4953
$ws.Range("K116").Value = 37055372
$ws.Range("L116").Value = 4953
$ws.Range("M116").Value = -37051930
$ws.Range("N116").Value = -11837

$ws.Range("H132").Value = 993775
$ws.Range("I132").Value = 1169006.9
$ws.Range("J132").Value = 30000
$ws.Range("K132").Value = 3507020.7
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -3504490.7
$ws.Range("N132").Value = -95060

$ws.Range("H137").Value = 2867.2812
$ws.Range("I137").Value = 1350.25
$ws.Range("J137").Value = 5395.6665
$ws.Range("K137").Value = 4050.75
$ws.Range("L137").Value = 16186.9995
$ws.Range("M137").Value = -1500.75
$ws.Range("N137").Value = -21286.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1196619.9
$ws.Range("I32").Value = 578.7083
$ws.Range("J32").Value = 5297332.5
$ws.Range("K32").Value = 578.7083
$ws.Range("L32").Value = 5297332.5
$ws.Range("M32").Value = -291.7083
$ws.Range("N32").Value = -5297906.5

$ws.Range("H74").Value = 2739.9412
$ws.Range("I74").Value = 1073.9
$ws.Range("J74").Value = 5120
$ws.Range("K74").Value = 1073.9
$ws.Range("L74").Value = 5120
$ws.Range("M74").Value = -199.9000000000001
$ws.Range("N74").Value = -6868

$ws.Range("H77").Value = 2739.9412
$ws.Range("I77").Value = 1073.9
$ws.Range("J77").Value = 5120
$ws.Range("K77").Value = 5369.5
$ws.Range("L77").Value = 25600
$ws.Range("M77").Value = -1001.5
$ws.Range("N77").Value = -34336

$ws.Range("H125").Value = 41999.668
$ws.Range("I125").Value = 40000
$ws.Range("J125").Value = 42999.5
$ws.Range("K125").Value = 40000
$ws.Range("L125").Value = 42999.5
$ws.Range("M125").Value = -35080
$ws.Range("N125").Value = -52839.5

$ws.Range("H132").Value = 1203366.2
$ws.Range("I132").Value = 2333938
$ws.Range("J132").Value = 72794.46000000001
$ws.Range("K132").Value = 7001814
$ws.Range("L132").Value = 218383.38
$ws.Range("M132").Value = -6999284
$ws.Range("N132").Value = -223443.38

$ws.Range("H133").Value = 153333
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 153333
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 153333
$ws.Range("N133").Value = -158393

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 26319148
$ws.Range("I105").Value = 90913800
$ws.Range("J105").Value = 2806.8147
$ws.Range("K105").Value = 90913800
$ws.Range("L105").Value = 2806.8147
$ws.Range("M105").Value = -90912053
$ws.Range("N105").Value = -6300.8147

$ws.Range("H107").Value = 4168117
$ws.Range("I107").Value = 4349209
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 4349209
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -4347289
$ws.Range("N107").Value = -6840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2631.9714
$ws.Range("I31").Value = 1884.5555
$ws.Range("J31").Value = 2890.6924
$ws.Range("K31").Value = 1884.5555
$ws.Range("L31").Value = 2890.6924
$ws.Range("M31").Value = -1589.5555
$ws.Range("N31").Value = -3480.6924

$ws.Range("H34").Value = 2631.9714
$ws.Range("I34").Value = 1884.5555
$ws.Range("J34").Value = 2890.6924
$ws.Range("K34").Value = 1884.5555
$ws.Range("L34").Value = 2890.6924
$ws.Range("M34").Value = -1682.5555
$ws.Range("N34").Value = -3294.6924

$ws.Range("H62").Value = 17799.666
$ws.Range("I62").Value = 17799.666
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 17799.666
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -17175.666

$ws.Range("H65").Value = 17799.666
$ws.Range("I65").Value = 17799.666
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 88998.33
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -85878.33

$ws.Range("H86").Value = 33302.75
$ws.Range("I86").Value = 32285.8
$ws.Range("J86").Value = 34997.668
$ws.Range("K86").Value = 32285.8
$ws.Range("L86").Value = 34997.668
$ws.Range("M86").Value = -31162.8
$ws.Range("N86").Value = -37243.668

$ws.Range("H89").Value = 33302.75
$ws.Range("I89").Value = 32285.8
$ws.Range("J89").Value = 34997.668
$ws.Range("K89").Value = 161429
$ws.Range("L89").Value = 174988.34
$ws.Range("M89").Value = -155813
$ws.Range("N89").Value = -186220.34

$ws.Range("H132").Value = 5394.8945
$ws.Range("I132").Value = 4044.2942
$ws.Range("J132").Value = 16875
$ws.Range("K132").Value = 12132.8826
$ws.Range("L132").Value = 50625
$ws.Range("M132").Value = -9602.882599999999
$ws.Range("N132").Value = -55685

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 308.5909
$ws.Range("I26").Value = 316.3158
$ws.Range("J26").Value = 259.66666
$ws.Range("K26").Value = 948.9474
$ws.Range("L26").Value = 778.9999799999999
$ws.Range("M26").Value = -660.9474
$ws.Range("N26").Value = -1354.99998

$ws.Range("H37").Value = 219995
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 219995
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 659985
$ws.Range("N37").Value = -660209

$ws.Range("H56").Value = 7487.5884
$ws.Range("I56").Value = 7487.5884
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 7487.5884
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -6957.5884

$ws.Range("H80").Value = 1997.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1997.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 5992.5
$ws.Range("N80").Value = -7864.5

$ws.Range("H83").Value = 1997.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1997.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17977.5
$ws.Range("N83").Value = -27337.5

$ws.Range("H127").Value = 4866.3335
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 4866.3335
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 14599.0005
$ws.Range("N127").Value = -24519.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 3870.7144
$ws.Range("I57").Value = 3870.7144
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 3870.7144
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -3050.7144
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 3225.2974
$ws.Range("I55").Value = 1997.0769
$ws.Range("J55").Value = 3890.5833
$ws.Range("K55").Value = 1997.0769
$ws.Range("L55").Value = 3890.5833
$ws.Range("M55").Value = -1824.0769
$ws.Range("N55").Value = -4236.5833

$ws.Range("H122").Value = 4168.76
$ws.Range("I122").Value = 3271.8462
$ws.Range("J122").Value = 5140.4165
$ws.Range("K122").Value = 9815.5386
$ws.Range("L122").Value = 15421.2495
$ws.Range("M122").Value = -7365.5386
$ws.Range("N122").Value = -20321.2495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2467.3948
$ws.Range("I96").Value = 1862.8
$ws.Range("J96").Value = 2861.6956
$ws.Range("K96").Value = 1862.8
$ws.Range("L96").Value = 2861.6956
$ws.Range("M96").Value = -489.8
$ws.Range("N96").Value = -5607.6956

$ws.Range("H106").Value = 49499.5
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 49499.5
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 49499.5
$ws.Range("N106").Value = -52023.5
